# Apply the "GUI Update (Add, Remove, and Edit working!)" edit to Sheet1 of the timesheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# C8 text "12:00pm-" becomes "12:00pm-5:15pm" (shared string edit in-place)
$ws.Range("C8").Value = "12:00pm-5:15pm"

# D8 (Hours) gets the numeric value 5.25
$ws.Range("D8").Value = 5.25

# E8 (Worked On) gets the new shared string describing the work done
$ws.Range("E8").Value = "Fix Commit Conflict, Implement Database Functionality to GUI"

# Update the active selection to E8, matching the saved view state
$ws.Range("E8").Select()
